# The document has four occurrences of a pattern that looks like:
#   <w:r>...<w:t>&lt;id&gt;</w:t></w:r>
#   <w:r>...<w:t>p041rXX</w:t></w:r>
#   <w:r>...<w:t>&lt;/id&gt;</w:t></w:r>
# Each triple of runs needs to collapse into a single run (keeping the
# formatting of the first run: the gold/Courier-New "tag" styling) whose
# text is the concatenation "<id>p041rXX</id>".
#
# Doing a Find.Execute with identical find/replace text over the full
# "<id>...</id>" span merges the matched runs into a single run that
# takes on the formatting of the first run of the match, which is
# exactly the transformation we want here - no visible text changes,
# only run (re)structuring.

$d = $word.ActiveDocument

$ids = @("p041r_1", "p041r_02", "p041r_03", "p041r_04")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $needle, 2) | Out-Null
}
